$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.395.14"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").Value = "1.873.88"
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("D4").Value = "'1.019"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.79%  "
$ws.Range("D5").Value = "'316.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").Value = "'1.017"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.99%  "
$ws.Range("D7").Value = "'0.5111"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.41%  "
$ws.Range("D8").Value = "'0.3941"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").Value = "'0.08437"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").Value = "'1.109"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.22%  "
$ws.Range("D11").Value = "'6.252"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").Value = "1.879.72"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("D13").Value = "'20.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "'7.234"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("D15").Value = "'1.019"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.25%  "
$ws.Range("D16").Value = "'0.00001109"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("D17").Value = "'91.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "'0.06769"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "'17.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("D20").Value = "'1.017"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").Value = "'5.948"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("D22").Value = "28.443.62"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").Value = "'11.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").Value = "'2.274"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").Value = "2.089.99"
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("D26").Value = "'161.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("D27").Value = "'20.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("D28").Value = "'2.359"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.29%  "
$ws.Range("D29").Value = "'126.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("D30").Value = "'0.1054"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "'1.043"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").Value = "'5.766"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.37%  "
$ws.Range("D33").Value = "'3.626"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.18%  "
$ws.Range("D34").Value = "'0.02428"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("D35").Value = "'0.06478"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.31%  "
$ws.Range("D36").Value = "'0.2172"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.60%  "
$ws.Range("D37").Value = "'8.851"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.71%  "
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("D39").Value = "'1.186"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").Value = "'0.6389"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.49%  "
$ws.Range("D41").Value = "'5.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("D42").Value = "'11.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.08%  "
$ws.Range("D43").Value = "'1.018"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").Value = "'0.6039"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").Value = "'12.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("D46").Value = "'3.715"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").Value = "'1.988"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.81%  "
$ws.Range("D48").Value = "'1.202"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.54%  "
$ws.Range("D49").Value = "'122.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").Value = "'1.204"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.15%  "
$ws.Range("D51").Value = "'0.06836"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.18%  "
